# Generate Report for Handoff
# Appends a new file-status row (aca95ef9-7219-4150-9599-1a05081b9a56) to the
# "Overview", "zh-cn" and "de-de" sheets of the localization-status report.

$wb = $excel.ActiveWorkbook

$fileId   = "aca95ef9-7219-4150-9599-1a05081b9a56"
$xlfHash  = "e0d0ee23577ace18383d04bf725f668e703be993"
$status   = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview" - row 9
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(9, 1).Value = "$fileId.md"
$wsOverview.Cells.Item(9, 2).Value = $status
$wsOverview.Cells.Item(9, 3).Value = $status
$wsOverview.Cells.Item(9, 4).Value = "2016-28-11 14:28:29"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsOverview.Range("A9").Font.Underline = $true
$wsOverview.Range("A9").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "zh-cn" - row 9
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Cells.Item(9, 1).Value = "$fileId.md"
$wsZhCn.Cells.Item(9, 2).Value = ".md"
$wsZhCn.Cells.Item(9, 3).Value = $status
$wsZhCn.Cells.Item(9, 4).Value = "$fileId.$xlfHash.zh-cn.xlf"
$wsZhCn.Cells.Item(9, 5).Value = "2016-03-11 14:28:26"
$wsZhCn.Cells.Item(9, 8).Value = "0001-01-01 00:00:00"
$wsZhCn.Cells.Item(9, 9).Value = "Include"

$wsZhCn.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/$fileId/e2e/$fileId.md",
    "",
    "",
    ".md"
)
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$xlfHash/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$fileId.$xlfHash.zh-cn.xlf",
    "",
    "",
    "$fileId.$xlfHash.zh-cn.xlf"
)

$wsZhCn.Range("A9").Font.Underline = $true
$wsZhCn.Range("A9").Font.Color = 15570276
$wsZhCn.Range("B9").Font.Underline = $true
$wsZhCn.Range("B9").Font.Color = 15570276
$wsZhCn.Range("D9").Font.Underline = $true
$wsZhCn.Range("D9").Font.Color = 15570276

# ---------------------------------------------------------------------
# Sheet "de-de" - row 9
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Cells.Item(9, 1).Value = "$fileId.md"
$wsDeDe.Cells.Item(9, 2).Value = ".md"
$wsDeDe.Cells.Item(9, 3).Value = $status
$wsDeDe.Cells.Item(9, 4).Value = "$fileId.$xlfHash.de-de.xlf"
$wsDeDe.Cells.Item(9, 5).Value = "2016-03-11 14:28:29"
$wsDeDe.Cells.Item(9, 8).Value = "0001-01-01 00:00:00"
$wsDeDe.Cells.Item(9, 9).Value = "Include"

$wsDeDe.Range("E9").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A9"),
    "https://github.com/OpenLocalizationTest/oltest/blob/$fileId/e2e/$fileId.md",
    "",
    "",
    "$fileId.md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B9"),
    "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/$fileId/e2e/$fileId.md",
    "",
    "",
    ".md"
)
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D9"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$xlfHash/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$fileId.$xlfHash.de-de.xlf",
    "",
    "",
    "$fileId.$xlfHash.de-de.xlf"
)

$wsDeDe.Range("A9").Font.Underline = $true
$wsDeDe.Range("A9").Font.Color = 15570276
$wsDeDe.Range("B9").Font.Underline = $true
$wsDeDe.Range("B9").Font.Color = 15570276
$wsDeDe.Range("D9").Font.Underline = $true
$wsDeDe.Range("D9").Font.Color = 15570276
